$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the JLCPCB Part#(optional) column (D) with the part numbers that were
# looked up for this BOM / pick-and-place pass. Values are prefixed with a
# leading apostrophe so they are stored as literal text (matching the existing
# quote-prefixed text style already used throughout column D).
$parts = @{
    2  = "C116352"
    3  = "C1779"
    4  = "C15850"
    5  = "C1711"
    6  = "C63339"
    7  = "C45783"
    8  = "C46653"
    9  = "C111566"
    10 = "C2297"
    11 = "C779674"
    12 = "C79382"
    14 = "C165948"
    15 = "C131339"
    16 = "C22438292"
    17 = "C154797"
    18 = "C218607"
    19 = "C2077400"
    20 = "C322246"
    21 = "C26023"
    22 = "C126687"
    23 = "C2828726"
    24 = "C25275"
    26 = "C17414"
    27 = "C25275"
    28 = "C17552"
    29 = "C26011"
    30 = "C25623"
    31 = "C967766"
    33 = "C54313"
    34 = "C86781"
    35 = "C3013946"
    36 = "C54313"
}

foreach ($row in $parts.Keys) {
    $ws.Range("D$row").Value = "'" + $parts[$row]
}

# Row 5 ended up re-wrapped (matches the wrapText cell style already used
# elsewhere in the sheet) when this value was entered.
$ws.Range("D5").WrapText = $true

# Row 26's designator wraps to two lines, so its height grew a bit once the
# sheet was touched again.
$ws.Rows.Item(26).RowHeight = 30

# Leave the view scrolled down near the bottom of the table, with the last
# filled-in cell selected, matching where editing finished.
$ws.Range("D37").Select()
